# Work: developed and tested InverseValidator, DirectValidator, PotenotValidator
#
# The "Journal" sheet gets a new journal entry. A blank row (row 27) is
# turned into a data row, and three new blank rows are inserted below it
# (before the existing "Сумма" block), shifting everything that used to be
# at row 28 onward down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Insert three new blank rows right after the (currently blank) row 27.
# This pushes the old rows 28.. (the Sum row, the second "web" block, the
# trailer rows, etc.) down to 31.. , matching the target layout.
$ws.Rows("28:30").Insert()

# The freshly inserted rows inherit odd auto styles from Insert(); restore
# the plain data-row look (same formatting as the rest of the blank rows in
# this table, e.g. row 26) by copying formats over them.
$ws.Range("A26:E26").Copy()
$ws.Range("A28:E30").PasteSpecial(-4122)

# Fill in the new journal entry in row 27.
$ws.Range("A27").Value2 = "InverseValidator, DirectValidator, PotenotValidator разработка и тесты"
$ws.Range("B27").Value2 = 45694
$ws.Range("C27").Value2 = 0.0972222222222222
$ws.Range("D27").Value2 = "Validator"

# Reflect the author's cursor position when they saved the file.
[void]$ws.Range("A5").Select()

Write-Host "Journal entry added"
